$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.992.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.622.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.620.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.092.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.865.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.625.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "338.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.98%  "
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0798"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "167.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0562"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.622"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0955"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
